# Update NATMI TPM output with newly recalculated edge-expression metrics.
# (Sending cluster / Ligand / Receptor / Target cluster labels are unchanged;
# only the computed numeric columns F, G, H, M-T are refreshed per row.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02551366666666667
$ws.Range("H2").Value = 0.076541
$ws.Range("M2").Value = 3.310473333333333
$ws.Range("N2").Value = 9.931419999999999
$ws.Range("O2").Value = 0.1683295705132556
$ws.Range("P2").Value = 0.1783110568845311
$ws.Range("Q2").Value = 0.08446231313555555
$ws.Range("R2").Value = 0.7601608182199999
$ws.Range("S2").Value = 0.1683295705132556
$ws.Range("T2").Value = 0.1783110568845311

# Row 3 (Target cluster: FAPs)
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02551366666666667
$ws.Range("H3").Value = 0.076541
$ws.Range("M3").Value = 3.966196333333333
$ws.Range("N3").Value = 11.898589
$ws.Range("O3").Value = 0.2016715007605908
$ws.Range("P3").Value = 0.2136300730433972
$ws.Range("Q3").Value = 0.1011922111832222
$ws.Range("R3").Value = 0.9107299006489999
$ws.Range("S3").Value = 0.2016715007605908
$ws.Range("T3").Value = 0.2136300730433972

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02551366666666667
$ws.Range("H4").Value = 0.076541
$ws.Range("M4").Value = 5.013056
$ws.Range("N4").Value = 15.039168
$ws.Range("O4").Value = 0.2549017854764673
$ws.Range("P4").Value = 0.2700167690767302
$ws.Range("Q4").Value = 0.1279014397653333
$ws.Range("R4").Value = 1.151112957888
$ws.Range("S4").Value = 0.2549017854764673
$ws.Range("T4").Value = 0.2700167690767302

# Row 5 (Target cluster: MuSCs)
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02551366666666667
$ws.Range("H5").Value = 0.076541
$ws.Range("M5").Value = 3.3026905
$ws.Range("N5").Value = 6.605381
$ws.Range("O5").Value = 0.1679338322424817
$ws.Range("P5").Value = 0.1185945682727144
$ws.Range("Q5").Value = 0.08426374452016666
$ws.Range("R5").Value = 0.505582467121
$ws.Range("S5").Value = 0.1679338322424817
$ws.Range("T5").Value = 0.1185945682727144

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02551366666666667
$ws.Range("H6").Value = 0.076541
$ws.Range("M6").Value = 4.074201666666667
$ws.Range("N6").Value = 12.222605
$ws.Range("O6").Value = 0.2071633110072045
$ws.Range("P6").Value = 0.2194475327226272
$ws.Range("Q6").Value = 0.1039478232561111
$ws.Range("R6").Value = 0.9355304093050001
$ws.Range("S6").Value = 0.2071633110072045
$ws.Range("T6").Value = 0.2194475327226272
